$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.202020202020202
$ws.Range("C2").Value = 0.5589225589225589
$ws.Range("J2").Value = 0.01683501683501683
$ws.Range("P2").Value = 0.1414141414141414
$ws.Range("S2").Value = 0.08080808080808081
$ws.Range("B3").Value = 0.002967359050445104
$ws.Range("C3").Value = 0.02967359050445104
$ws.Range("J3").Value = 0.03560830860534125
$ws.Range("P3").Value = 0.7507418397626113
$ws.Range("S3").Value = 0.1810089020771513
$ws.Range("J4").Value = 0.05504587155963303
$ws.Range("P4").Value = 0.7064220183486238
$ws.Range("S4").Value = 0.2385321100917431
$ws.Range("B6").Value = 0.03485838779956427
$ws.Range("D6").Value = 0.0130718954248366
$ws.Range("E6").Value = 0.002178649237472767
$ws.Range("F6").Value = 0.05882352941176471
$ws.Range("J6").Value = 0.2549019607843137
$ws.Range("O6").Value = 0.01742919389978214
$ws.Range("Q6").Value = 0.159041394335512
$ws.Range("R6").Value = 0.08714596949891068
$ws.Range("S6").Value = 0.3725490196078431
$ws.Range("B7").Value = 0.09848484848484848
$ws.Range("D7").Value = 0.02525252525252525
$ws.Range("E7").Value = 0.005050505050505051
$ws.Range("F7").Value = 0.04545454545454546
$ws.Range("J7").Value = 0.1363636363636364
$ws.Range("O7").Value = 0.01767676767676768
$ws.Range("Q7").Value = 0.154040404040404
$ws.Range("R7").Value = 0.06818181818181818
$ws.Range("S7").Value = 0.4494949494949495
$ws.Range("B8").Value = 0.09808102345415778
$ws.Range("D8").Value = 0.01492537313432836
$ws.Range("E8").Value = 0.001066098081023454
$ws.Range("F8").Value = 0.06076759061833688
$ws.Range("J8").Value = 0.1204690831556503
$ws.Range("O8").Value = 0.0138592750533049
$ws.Range("Q8").Value = 0.185501066098081
$ws.Range("R8").Value = 0.1087420042643923
$ws.Range("S8").Value = 0.3965884861407249
$ws.Range("B9").Value = 0.08631578947368421
$ws.Range("D9").Value = 0.01894736842105263
$ws.Range("F9").Value = 0.06736842105263158
$ws.Range("J9").Value = 0.1136842105263158
$ws.Range("O9").Value = 0.01263157894736842
$ws.Range("Q9").Value = 0.1852631578947369
$ws.Range("R9").Value = 0.08842105263157894
$ws.Range("S9").Value = 0.4273684210526316
$ws.Range("B10").Value = 0.09659485753995831
$ws.Range("D10").Value = 0.02571230020847811
$ws.Range("E10").Value = 0.0003474635163307853
$ws.Range("F10").Value = 0.06497567755385684
$ws.Range("J10").Value = 0.1191799861014593
$ws.Range("O10").Value = 0.01667824878387769
$ws.Range("Q10").Value = 0.2282835302293259
$ws.Range("R10").Value = 0.08895066018068103
$ws.Range("S10").Value = 0.359277275886032
$ws.Range("F11").Value = 0.001615508885298869
$ws.Range("G11").Value = 0.148626817447496
$ws.Range("J11").Value = 0.09531502423263329
$ws.Range("K11").Value = 0.1987075928917609
$ws.Range("L11").Value = 0.5444264943457189
$ws.Range("S11").Value = 0.01130856219709208
$ws.Range("G12").Value = 0.7225274725274725
$ws.Range("J12").Value = 0.1868131868131868
$ws.Range("K12").Value = 0.01648351648351648
$ws.Range("L12").Value = 0.06593406593406594
$ws.Range("S12").Value = 0.008241758241758242
$ws.Range("F15").Value = 0.01965065502183406
$ws.Range("H15").Value = 0.1550218340611354
$ws.Range("I15").Value = 0.07860262008733625
$ws.Range("J15").Value = 0.3558951965065502
$ws.Range("K15").Value = 0.07205240174672489
$ws.Range("M15").Value = 0.002183406113537118
$ws.Range("O15").Value = 0.06550218340611354
$ws.Range("S15").Value = 0.2510917030567685
$ws.Range("F16").Value = 0.01732673267326733
$ws.Range("H16").Value = 0.1905940594059406
$ws.Range("I16").Value = 0.09900990099009901
$ws.Range("J16").Value = 0.3836633663366337
$ws.Range("K16").Value = 0.1089108910891089
$ws.Range("M16").Value = 0.0198019801980198
$ws.Range("N16").Value = 0.002475247524752475
$ws.Range("O16").Value = 0.05445544554455446
$ws.Range("S16").Value = 0.1237623762376238
$ws.Range("F17").Value = 0.01620591039084843
$ws.Range("H17").Value = 0.1725452812202097
$ws.Range("I17").Value = 0.1010486177311725
$ws.Range("J17").Value = 0.4270734032411821
$ws.Range("K17").Value = 0.0905624404194471
$ws.Range("M17").Value = 0.01715919923736892
$ws.Range("N17").Value = 0.0009532888465204957
$ws.Range("O17").Value = 0.05910390848427073
$ws.Range("S17").Value = 0.11534795042898
$ws.Range("F18").Value = 0.01518438177874186
$ws.Range("H18").Value = 0.158351409978308
$ws.Range("I18").Value = 0.09327548806941431
$ws.Range("J18").Value = 0.4598698481561822
$ws.Range("K18").Value = 0.09327548806941431
$ws.Range("M18").Value = 0.008676789587852495
$ws.Range("N18").Value = 0.002169197396963124
$ws.Range("O18").Value = 0.06290672451193059
$ws.Range("S18").Value = 0.1062906724511931
$ws.Range("F19").Value = 0.01406365655070318
$ws.Range("H19").Value = 0.1991117690599556
$ws.Range("I19").Value = 0.0921539600296077
$ws.Range("J19").Value = 0.3963730569948187
$ws.Range("K19").Value = 0.09844559585492228
$ws.Range("M19").Value = 0.01739452257586973
$ws.Range("N19").Value = 0.002590673575129534
$ws.Range("O19").Value = 0.06809770540340489
$ws.Range("S19").Value = 0.1117690599555885
